$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 538, pushing existing rows 538-581 down to 540-583
$ws.Range("538:539").Insert()

# Fill in new row 538 with fresh data
$ws.Cells.Item(538, 1).Value2 = 10
$ws.Cells.Item(538, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(538, 3).Value2 = "La Araucanía"
$ws.Cells.Item(538, 4).Value2 = 44746
$ws.Cells.Item(538, 5).Value2 = 9
$ws.Cells.Item(538, 6).Value2 = 100112045
$ws.Cells.Item(538, 7).Value2 = "Zapallo"
$ws.Cells.Item(538, 8).Value2 = "Camote"
$ws.Cells.Item(538, 9).Value2 = "1a (guarda)"
$ws.Cells.Item(538, 10).Value2 = 450
$ws.Cells.Item(538, 11).Value2 = 800
$ws.Cells.Item(538, 12).Value2 = 800
$ws.Cells.Item(538, 13).Value2 = 800
$ws.Cells.Item(538, 14).Value2 = "$/kilo (volumen en unidades)"
$ws.Cells.Item(538, 15).Value2 = "Región del Maule"
$ws.Cells.Item(538, 16).Value2 = 800
$ws.Cells.Item(538, 17).Value2 = 1
$ws.Cells.Item(538, 18).Value2 = "Hortaliza"

# Fill in new row 539 with fresh data
$ws.Cells.Item(539, 1).Value2 = 10
$ws.Cells.Item(539, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(539, 3).Value2 = "La Araucanía"
$ws.Cells.Item(539, 4).Value2 = 44746
$ws.Cells.Item(539, 5).Value2 = 9
$ws.Cells.Item(539, 6).Value2 = 100112045
$ws.Cells.Item(539, 7).Value2 = "Zapallo"
$ws.Cells.Item(539, 8).Value2 = "Camote"
$ws.Cells.Item(539, 9).Value2 = "2a (guarda)"
$ws.Cells.Item(539, 10).Value2 = 890
$ws.Cells.Item(539, 11).Value2 = 500
$ws.Cells.Item(539, 12).Value2 = 500
$ws.Cells.Item(539, 13).Value2 = 500
$ws.Cells.Item(539, 14).Value2 = "$/kilo (volumen en unidades)"
$ws.Cells.Item(539, 15).Value2 = "Región del Maule"
$ws.Cells.Item(539, 16).Value2 = 500
$ws.Cells.Item(539, 17).Value2 = 1
$ws.Cells.Item(539, 18).Value2 = "Hortaliza"
